$wb = $excel.ActiveWorkbook

# The Czech sheet keeps its own row selection (no longer the active/visible tab)
# after the new sheet becomes active - update its remembered selection first.
$czech = $wb.Worksheets.Item("Czech")
$czech.Range("B19").Select() | Out-Null

# Duplicate the Czech sheet (same layout/styles/merges) right after itself, then
# turn the copy into the new "Swiss" market sheet.
$czech.Copy($null, $czech)
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# Update the market-specific text.
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2642"

# Make the new sheet the active tab with its own remembered selection.
$swiss.Activate() | Out-Null
$swiss.Range("B7").Select() | Out-Null
